$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Replace the paragraph's visible text (everything except the trailing
# paragraph mark) with the new placeholder token. This collapses the
# two existing runs (the id token + trailing space run) into a single
# run, inheriting the first run's character formatting.
$r = $p1.Range
$r.End = $r.End - 1
$r.Text = "**ID__AFFARS_SUBPART_5312_4__ID**"

# Add a paragraph border (top/left/bottom/right, 5pt space, default
# line) around the first paragraph.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# Widen the left indent from 120 to 225 twips (120/20=6pt -> 225/20=11.25pt).
$p1.Format.LeftIndent = 11.25
